$word.UserName = "Rob"
$d = $word.ActiveDocument
$d.TrackRevisions = $true

# -----------------------------------------------------------------
# 1) "Be able to withstand exposure to rain" -> delete "Be able to "
#    leaving "withstand exposure to rain", then append a new tracked
#    insertion " (Excluding payload)" after the existing
#    " while powered down" insertion.
# -----------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Be able to withstand exposure to rain", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r1.Find.Found) {
  $delRange = $d.Range($r1.Start, $r1.Start + 11)
  $delRange.Delete()
}

$r2 = $d.Content
$r2.Find.Execute("while powered down", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r2.Find.Found) {
  $insPos = $r2.End
  $ip2 = $d.Range($insPos, $insPos)
  $ip2.InsertBefore(" (Excluding payload)")
}

# -----------------------------------------------------------------
# 2) "Allow manual control of manual axes" -> "Allow manual control
#    of individual axes" (the whole phrase is already a pending
#    insertion, so editing inside it just shrinks/grows that
#    insertion rather than creating new del/ins pairs).
# -----------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("manual axes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r3.Find.Found) {
  $wordRange = $d.Range($r3.Start, $r3.Start + 6)
  $insPos3 = $wordRange.Start
  $wordRange.Delete()
  $ip3 = $d.Range($insPos3, $insPos3)
  $ip3.InsertBefore("individual")
}

# -----------------------------------------------------------------
# 3) New paragraph after "Support remote operation":
#    "Communicate with other devices via UDP"
# -----------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Support remote operation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r4.Find.Found) {
  $r4.Collapse(0)
  $r4.InsertAfter([char]13)

  $r5 = $d.Content
  $r5.Find.Execute("Support remote operation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
  $r5.Collapse(0)
  $newParaStart = $r5.Start + 1
  $ip5 = $d.Range($newParaStart, $newParaStart)
  $ip5.InsertBefore("Communicate with other devices via UDP")
}
